$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (row 1) ---
$ws.Range("B1").Value = '$ bold(''All'')'
$ws.Range("C1").Value = '$ bold(''Europe'')'
$ws.Range("D1").Value = 'France'
$ws.Range("E1").Value = 'Germany'
$ws.Range("F1").Value = 'Italy'
$ws.Range("G1").Value = 'Poland'
$ws.Range("H1").Value = 'Spain'
$ws.Range("I1").Value = 'United Kingdom'
$ws.Range("J1").Value = 'Switzerland'
$ws.Range("K1").Value = 'Japan'
$ws.Range("L1").Value = 'Russia'
$ws.Range("M1").Value = 'Saudi Arabia'
$ws.Range("N1").Value = 'USA'

# --- Row labels (column A, rows 2-25) ---
$ws.Range("A2").Value = 'Money; own income; cost of living; inflation'
$ws.Range("A3").Value = 'Health; healthcare system'
$ws.Range("A4").Value = 'Own country referred'
$ws.Range("A5").Value = 'Family; children; childcare'
$ws.Range("A6").Value = 'War; peace'
$ws.Range("A7").Value = 'Work; (un)employment; business'
$ws.Range("A8").Value = 'Nothing; don''t know; empty'
$ws.Range("A9").Value = 'Economy'
$ws.Range("A10").Value = 'Government; president'
$ws.Range("A11").Value = 'International issues'
$ws.Range("A12").Value = 'Inflation; cost of living'
$ws.Range("A13").Value = 'Poverty; inequality'
$ws.Range("A14").Value = 'Tax system; welfare benefits; public services'
$ws.Range("A15").Value = 'Old age; retirement; ageing society'
$ws.Range("A16").Value = 'Criticism of immigration; national preference'
$ws.Range("A17").Value = 'Housing'
$ws.Range("A18").Value = 'Security; violence; crime; judicial system'
$ws.Range("A19").Value = 'Criticism of far right; Trump; tariffs'
$ws.Range("A20").Value = 'Environment; climate change'
$ws.Range("A21").Value = 'Rights; democracy; freedom; slavery'
$ws.Range("A22").Value = 'Discrimination; gender inequality; racism; LGBT'
$ws.Range("A23").Value = 'Happiness; peace of mind'
$ws.Range("A24").Value = 'Trump'
$ws.Range("A25").Value = 'Relationships; love; emotions'

# --- Data matrix (columns B-N, rows 2-25) ---
$ws.Range("B2").Value = 0.182321037798967
$ws.Range("C2").Value = 0.161772357883732
$ws.Range("D2").Value = 0.151628165867824
$ws.Range("E2").Value = 0.15260419576115
$ws.Range("F2").Value = 0.164293045891492
$ws.Range("G2").Value = 0.167992838290846
$ws.Range("H2").Value = 0.141646939023308
$ws.Range("I2").Value = 0.193962658106887
$ws.Range("J2").Value = 0.151337475119412
$ws.Range("K2").Value = 0.205749694909397
$ws.Range("L2").Value = 0.128357235867269
$ws.Range("M2").Value = 0.160511571395735
$ws.Range("N2").Value = 0.22136380470097
$ws.Range("B3").Value = 0.100995060768002
$ws.Range("C3").Value = 0.126903713743287
$ws.Range("D3").Value = 0.101500215876694
$ws.Range("E3").Value = 0.107442105698588
$ws.Range("F3").Value = 0.14622572265218
$ws.Range("G3").Value = 0.123299614530166
$ws.Range("H3").Value = 0.159479355203936
$ws.Range("I3").Value = 0.140689286612206
$ws.Range("J3").Value = 0.103994548882937
$ws.Range("K3").Value = 0.0438097949544678
$ws.Range("L3").Value = 0.0918208304430512
$ws.Range("M3").Value = 0.0481397016758384
$ws.Range("N3").Value = 0.102648655642718
$ws.Range("B4").Value = 0.088154808321385
$ws.Range("C4").Value = 0.0907163404224399
$ws.Range("D4").Value = 0.106882912115324
$ws.Range("E4").Value = 0.0947465259960133
$ws.Range("F4").Value = 0.0632037251264596
$ws.Range("G4").Value = 0.110593408120899
$ws.Range("H4").Value = 0.0788491814450573
$ws.Range("I4").Value = 0.0959198410759693
$ws.Range("J4").Value = 0.0600154515657946
$ws.Range("K4").Value = 0.0838073240940032
$ws.Range("L4").Value = 0.0608786217043796
$ws.Range("M4").Value = 0.0602178270789178
$ws.Range("N4").Value = 0.101040837559806
$ws.Range("B5").Value = 0.0711807113712516
$ws.Range("C5").Value = 0.0651884757861758
$ws.Range("D5").Value = 0.0536855985040587
$ws.Range("E5").Value = 0.042692829320797
$ws.Range("F5").Value = 0.0689439875428471
$ws.Range("G5").Value = 0.0767926429993233
$ws.Range("H5").Value = 0.0529653279236066
$ws.Range("I5").Value = 0.10728342826472
$ws.Range("J5").Value = 0.0428373737667234
$ws.Range("K5").Value = 0.0731479352716921
$ws.Range("L5").Value = 0.0662071089146564
$ws.Range("M5").Value = 0.102037616300348
$ws.Range("N5").Value = 0.0766206382747447
$ws.Range("B6").Value = 0.060736531600857
$ws.Range("C6").Value = 0.0942459172240121
$ws.Range("D6").Value = 0.0676862447367624
$ws.Range("E6").Value = 0.110163507125438
$ws.Range("F6").Value = 0.140402778999085
$ws.Range("G6").Value = 0.13974013453665
$ws.Range("H6").Value = 0.0649441752978565
$ws.Range("I6").Value = 0.0542563589055744
$ws.Range("J6").Value = 0.0856719423556321
$ws.Range("K6").Value = 0.0375112362668072
$ws.Range("L6").Value = 0.0330590623121272
$ws.Range("M6").Value = 0.0420394789499737
$ws.Range("N6").Value = 0.0452396364869285
$ws.Range("B7").Value = 0.0553504218200639
$ws.Range("C7").Value = 0.0607223120240159
$ws.Range("D7").Value = 0.0719609983384065
$ws.Range("E7").Value = 0.0501311815999043
$ws.Range("F7").Value = 0.0822118045049334
$ws.Range("G7").Value = 0.0399942430095793
$ws.Range("H7").Value = 0.0679040031831562
$ws.Range("I7").Value = 0.0538759262784276
$ws.Range("J7").Value = 0.0316915930729731
$ws.Range("K7").Value = 0.0426050988618772
$ws.Range("L7").Value = 0.0534989239387084
$ws.Range("M7").Value = 0.101771032271278
$ws.Range("N7").Value = 0.050970650353524
$ws.Range("B8").Value = 0.0469900026163117
$ws.Range("C8").Value = 0.0409974052476144
$ws.Range("D8").Value = 0.0542745978974009
$ws.Range("E8").Value = 0.0432300470705285
$ws.Range("F8").Value = 0.029744754106229
$ws.Range("G8").Value = 0.0663865396156052
$ws.Range("H8").Value = 0.0335897645052565
$ws.Range("I8").Value = 0.0279930616552985
$ws.Range("J8").Value = 0.0247185603776682
$ws.Range("K8").Value = 0.134585103893787
$ws.Range("L8").Value = 0.0426251947414407
$ws.Range("M8").Value = 0.0358156854121428
$ws.Range("N8").Value = 0.0221130651051977
$ws.Range("B9").Value = 0.0428925022529864
$ws.Range("C9").Value = 0.0387734005589055
$ws.Range("D9").Value = 0.011474957650423
$ws.Range("E9").Value = 0.0437714005759224
$ws.Range("F9").Value = 0.0626502583920119
$ws.Range("G9").Value = 0.0110333227815035
$ws.Range("H9").Value = 0.0514157193779358
$ws.Range("I9").Value = 0.0457402413085588
$ws.Range("J9").Value = 0.0309094204739583
$ws.Range("K9").Value = 0.035094409747487
$ws.Range("L9").Value = 0.0134987979493771
$ws.Range("M9").Value = 0.0199711849845625
$ws.Range("N9").Value = 0.0651946587957189
$ws.Range("B10").Value = 0.0403203062652297
$ws.Range("C10").Value = 0.0285615312623572
$ws.Range("D10").Value = 0.0320049726035638
$ws.Range("E10").Value = 0.0241596296022602
$ws.Range("F10").Value = 0.019577662110769
$ws.Range("G10").Value = 0.0180001922636902
$ws.Range("H10").Value = 0.0270399273282725
$ws.Range("I10").Value = 0.0475792780892838
$ws.Range("J10").Value = 0.0189273882550392
$ws.Range("K10").Value = 0.0385965646907139
$ws.Range("L10").Value = 0.00223387876409278
$ws.Range("M10").Value = 0.000862233596311687
$ws.Range("N10").Value = 0.074073346356152
$ws.Range("B11").Value = 0.0389960678335084
$ws.Range("C11").Value = 0.050449008233722
$ws.Range("D11").Value = 0.0547130842188104
$ws.Range("E11").Value = 0.0590899980535088
$ws.Range("F11").Value = 0.0512751354655843
$ws.Range("G11").Value = 0.0539464981958894
$ws.Range("H11").Value = 0.0394128679957092
$ws.Range("I11").Value = 0.0413631073311336
$ws.Range("J11").Value = 0.0435924668129858
$ws.Range("K11").Value = 0.0304564891017079
$ws.Range("L11").Value = 0.017840134031664
$ws.Range("M11").Value = 0.0777988628197214
$ws.Range("N11").Value = 0.0349429275482588
$ws.Range("B12").Value = 0.0384638187662071
$ws.Range("C12").Value = 0.0364601289385647
$ws.Range("D12").Value = 0.0150595822264171
$ws.Range("E12").Value = 0.0227496674399152
$ws.Range("F12").Value = 0.0489918617466034
$ws.Range("G12").Value = 0.0284025769810721
$ws.Range("H12").Value = 0.0286133153886044
$ws.Range("I12").Value = 0.0768024795846809
$ws.Range("J12").Value = 0.018779014651579
$ws.Range("K12").Value = 0.0234399934012472
$ws.Range("L12").Value = 0.0107231605115973
$ws.Range("M12").Value = 0.00642657788541434
$ws.Range("N12").Value = 0.0613077559670346
$ws.Range("B13").Value = 0.0373562141560839
$ws.Range("C13").Value = 0.0603739511218419
$ws.Range("D13").Value = 0.0538315457227264
$ws.Range("E13").Value = 0.0700313444742106
$ws.Range("F13").Value = 0.0636993820136352
$ws.Range("G13").Value = 0.0677494356427138
$ws.Range("H13").Value = 0.0590225495655187
$ws.Range("I13").Value = 0.051246217617998
$ws.Range("J13").Value = 0.0378364707404913
$ws.Range("K13").Value = 0.0312099154110523
$ws.Range("L13").Value = 0.021268099021952
$ws.Range("M13").Value = 0.0303083114789396
$ws.Range("N13").Value = 0.0210887233519309
$ws.Range("B14").Value = 0.0336191355177953
$ws.Range("C14").Value = 0.0292223231542092
$ws.Range("D14").Value = 0.0221765212170327
$ws.Range("E14").Value = 0.0416974809332313
$ws.Range("F14").Value = 0.0285423088877675
$ws.Range("G14").Value = 0.0176825445547167
$ws.Range("H14").Value = 0.0228942615202653
$ws.Range("I14").Value = 0.0337659725984563
$ws.Range("J14").Value = 0.0187933388549178
$ws.Range("K14").Value = 0.103938267948727
$ws.Range("L14").Value = 0.00138736039927219
$ws.Range("M14").Value = 0.00208813427186945
$ws.Range("N14").Value = 0.0274075262990083
$ws.Range("B15").Value = 0.0332160574775112
$ws.Range("C15").Value = 0.032541579467393
$ws.Range("D15").Value = 0.0273958857117257
$ws.Range("E15").Value = 0.0684305271844509
$ws.Range("F15").Value = 0.0124141813050852
$ws.Range("G15").Value = 0.0153742552149
$ws.Range("H15").Value = 0.0186892581057893
$ws.Range("I15").Value = 0.0325524793108507
$ws.Range("J15").Value = 0.0177829722154919
$ws.Range("K15").Value = 0.0590920345710629
$ws.Range("L15").Value = 0.0455635490148159
$ws.Range("M15").Value = 0.0039835499686648
$ws.Range("N15").Value = 0.0210925808322273
$ws.Range("B16").Value = 0.0319649704467399
$ws.Range("C16").Value = 0.0529305712288578
$ws.Range("D16").Value = 0.0357057845818555
$ws.Range("E16").Value = 0.0780398179127967
$ws.Range("F16").Value = 0.0274281677578105
$ws.Range("G16").Value = 0.0350434292800761
$ws.Range("H16").Value = 0.0321540740533865
$ws.Range("I16").Value = 0.0874578351856686
$ws.Range("J16").Value = 0.0521707440391609
$ws.Range("K16").Value = 0.00656675605998006
$ws.Range("L16").Value = 0.00476136709085499
$ws.Range("M16").Value = 0.00199639412499814
$ws.Range("N16").Value = 0.0323753859474656
$ws.Range("B17").Value = 0.0304488376315802
$ws.Range("C17").Value = 0.0269735173134916
$ws.Range("D17").Value = 0.0213216678839544
$ws.Range("E17").Value = 0.0154592879988818
$ws.Range("F17").Value = 0.0229261585613421
$ws.Range("G17").Value = 0.0330329698045628
$ws.Range("H17").Value = 0.0560473240220591
$ws.Range("I17").Value = 0.0278306333783973
$ws.Range("J17").Value = 0.0144175784657228
$ws.Range("K17").Value = 0.0146534566511689
$ws.Range("L17").Value = 0.0592517001731302
$ws.Range("M17").Value = 0.0319695272437633
$ws.Range("N17").Value = 0.0282176459782966
$ws.Range("B18").Value = 0.0302795812721631
$ws.Range("C18").Value = 0.0281598134127114
$ws.Range("D18").Value = 0.0163917340461058
$ws.Range("E18").Value = 0.020089009677194
$ws.Range("F18").Value = 0.0395814308851999
$ws.Range("G18").Value = 0.0129502889634224
$ws.Range("H18").Value = 0.023924431046729
$ws.Range("I18").Value = 0.0531799879492298
$ws.Range("J18").Value = 0.0176989396518012
$ws.Range("K18").Value = 0.0217650060909604
$ws.Range("L18").Value = 0.00586038734211524
$ws.Range("M18").Value = 0.0111451584411178
$ws.Range("N18").Value = 0.0481378270190486
$ws.Range("B19").Value = 0.029472401622532
$ws.Range("C19").Value = 0.0157395055510232
$ws.Range("D19").Value = 0.0176020083086458
$ws.Range("E19").Value = 0.0233644205389716
$ws.Range("F19").Value = 0.013900924047365
$ws.Range("G19").Value = 0.00839006350112947
$ws.Range("H19").Value = 0.00949857755645115
$ws.Range("I19").Value = 0.0149198636575293
$ws.Range("J19").Value = 0.0140747856649723
$ws.Range("K19").Value = 0.026397328866856
$ws.Range("L19").Value = 0.00140206063911074
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0.0608546638483413
$ws.Range("B20").Value = 0.0290290271298219
$ws.Range("C20").Value = 0.0437536500173379
$ws.Range("D20").Value = 0.0235830505733939
$ws.Range("E20").Value = 0.0548185750778946
$ws.Range("F20").Value = 0.0698277052752728
$ws.Range("G20").Value = 0.0172512797273797
$ws.Range("H20").Value = 0.0355182686179697
$ws.Range("I20").Value = 0.0469875691317198
$ws.Range("J20").Value = 0.0467077004351659
$ws.Range("K20").Value = 0.0141680897206946
$ws.Range("L20").Value = 0.0016522224901023
$ws.Range("M20").Value = 0.0497244153114467
$ws.Range("N20").Value = 0.0279939370302207
$ws.Range("B21").Value = 0.0273840558766951
$ws.Range("C21").Value = 0.0185734738812219
$ws.Range("D21").Value = 0.00993376996509722
$ws.Range("E21").Value = 0.0219503872991804
$ws.Range("F21").Value = 0.0190755572980811
$ws.Range("G21").Value = 0.0116551769289037
$ws.Range("H21").Value = 0.0191475138528719
$ws.Range("I21").Value = 0.026083694384323
$ws.Range("J21").Value = 0.0185893481653607
$ws.Range("K21").Value = 0.00855123830249976
$ws.Range("L21").Value = 0.00643806063116806
$ws.Range("M21").Value = 0.0424901047589513
$ws.Range("N21").Value = 0.0524220495695358
$ws.Range("B22").Value = 0.0247509642480483
$ws.Range("C22").Value = 0.0235491585698141
$ws.Range("D22").Value = 0.0247295118635322
$ws.Range("E22").Value = 0.0199228701141335
$ws.Range("F22").Value = 0.0239889324483805
$ws.Range("G22").Value = 0.00842241494155756
$ws.Range("H22").Value = 0.0215439821623808
$ws.Range("I22").Value = 0.0373860187913463
$ws.Range("J22").Value = 0.0199103360766923
$ws.Range("K22").Value = 0.0182662703177434
$ws.Range("L22").Value = 0.00159299207370676
$ws.Range("M22").Value = 0.0270490585304728
$ws.Range("N22").Value = 0.0383410720026572
$ws.Range("B23").Value = 0.0223801820434425
$ws.Range("C23").Value = 0.0278669994830184
$ws.Range("D23").Value = 0.0229189212280595
$ws.Range("E23").Value = 0.0180669488001356
$ws.Range("F23").Value = 0.043761530187622
$ws.Range("G23").Value = 0.00404139647620237
$ws.Range("H23").Value = 0.0227917901718656
$ws.Range("I23").Value = 0.0506642953767923
$ws.Range("J23").Value = 0.0107545909645981
$ws.Range("K23").Value = 0.00815618955253707
$ws.Range("L23").Value = 0.0105472628002141
$ws.Range("M23").Value = 0.0125275022226377
$ws.Range("N23").Value = 0.0276605926752763
$ws.Range("B24").Value = 0.0219291560180485
$ws.Range("C24").Value = 0.00961126540544619
$ws.Range("D24").Value = 0.0149985738479522
$ws.Range("E24").Value = 0.0114707836995252
$ws.Range("F24").Value = 0.00722338775958125
$ws.Range("G24").Value = 0.00306882833325286
$ws.Range("H24").Value = 0.000876257471860367
$ws.Range("I24").Value = 0.0136748148581175
$ws.Range("J24").Value = 0.0140747856649723
$ws.Range("K24").Value = 0.0209611471772189
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0.0475907360468383
$ws.Range("B25").Value = 0.0194397523322192
$ws.Range("C25").Value = 0.018120854570026
$ws.Range("D25").Value = 0.0225301695013904
$ws.Range("E25").Value = 0.0231518048949187
$ws.Range("F25").Value = 0.01332380993786
$ws.Range("G25").Value = 0.0107414062601043
$ws.Range("H25").Value = 0.0189943771325983
$ws.Range("I25").Value = 0.0168612567712371
$ws.Range("J25").Value = 0.00674796494606546
$ws.Range("K25").Value = 0.00398912595078207
$ws.Range("L25").Value = 0.0194287855799425
$ws.Range("M25").Value = 0.0153588440444347
$ws.Range("N25").Value = 0.0273943271811495
